# The sheet lists (group, label) pairs produced by a compact-letter-display
# script. In this revision, the "CSS x Reduced" and "Grassland x Reduced"
# rows were reordered (A2 <-> A3), while everything else on the sheet stays
# the same. The leftover reviewer comment that had been attached to B5 is
# also cleaned up as part of this pass.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap the group names in A2 and A3.
$ws.Range("A2").Value = "Grassland x Reduced"
$ws.Range("A3").Value = "CSS x Reduced"

# Remove the stale threaded comment that was left on B5.
$commentCell = $ws.Range("B5")
if ($commentCell.Comment -ne $null) {
    $commentCell.Comment.Delete()
}
